$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has a "VendorList" header in row 1 (bold Calibri,
# style used by A1) followed by the real vendor rows (2-21) in a plain
# Consolas style. The header row is being dropped, the vendor rows shift
# up by one, and every remaining row adopts the header's (bold) style.

# 1) Stamp the header's formatting onto the vendor rows before the header
#    row disappears, so the shift-up ends with a single, shared style.
$ws.Range("A1").Copy()
$ws.Range("A2:A21").PasteSpecial(-4122)

# 2) Delete the header row outright; rows 2-21 (Facebook ... Zomato)
#    shift up to become rows 1-20.
$ws.Rows.Item(1).Delete()

# 3) Leave the selection on the new last row, A20, matching the saved
#    cursor position.
[void]$ws.Range("A20").Select()
